$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Copy the formatting (styles/borders/fill) of row 2's D:H range down to row 3
# so the new cells pick up the same visual style used by the existing table.
$ws.Range("D2:H2").Copy()
$ws.Range("D3:H3").PasteSpecial(-4122)

# Fill in the new "remove bus trans" form row values
$ws.Range("D3").Value = 42
$ws.Range("E3").Value = "Selecione um item"
$ws.Range("F3").Value = "ComboBox 9"
$ws.Range("G3").Value = 15

# Match the author's final cell selection recorded in the sheet view
$ws.Range("D3").Select()
